$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.157.98"
$ws.Range("E2").Value = "  -5.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.564.02"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.33"
$ws.Range("E5").Value = "  -5.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.05"
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.560.39"
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -8.52%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.695"
$ws.Range("E10").Value = "  -9.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  -15.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000351"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.47"
$ws.Range("E13").Value = "  -7.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.111.11"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.35"
$ws.Range("E15").Value = "  -5.62%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.555.39"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.86"
$ws.Range("E18").Value = "  +3.84%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.95"
$ws.Range("E19").Value = "  -7.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "64.220.09"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("E21").Value = "  -9.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.69"
$ws.Range("E22").Value = "  -13.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.05"
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.26"
$ws.Range("E24").Value = "  -7.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.91"
$ws.Range("E25").Value = "  -7.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.47"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.33"
$ws.Range("E27").Value = "  +9.40%  "
$ws.Range("E28").Value = "  -9.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.92"
$ws.Range("E29").Value = "  -11.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.00"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.89"
$ws.Range("E33").Value = "  -7.08%  "
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.69"
$ws.Range("E35").Value = "  -7.59%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.32"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0443"
$ws.Range("E38").Value = "  -10.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +19.93%  "
$ws.Range("E41").Value = "  -9.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0625"
$ws.Range("E42").Value = "  -10.87%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.08"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.05"
$ws.Range("E44").Value = "  +13.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.11"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("E47").Value = "  -9.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.02"
$ws.Range("E49").Value = "  +15.99%  "
$ws.Range("E50").Value = "  -7.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.279"
$ws.Range("E51").Value = "  -9.54%  "
